$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New full data set (rows 2-18), header row (row 1) unchanged
$data = @(
    @("Keyonte George", "PG,SG", "Utah Jazz"),
    @("Dalano Banton", "SG,SF", "Portland Trail Blazers"),
    @("Tobias Harris", "SF,PF", "Detroit Pistons"),
    @("Lauri Markkanen", "SF,PF", "Utah Jazz"),
    @("Christian Braun", "SG,SF", "Denver Nuggets"),
    @("Zach LaVine", "SG,SF", "Chicago Bulls"),
    @("John Collins", "PF,C", "Utah Jazz"),
    @("Jalen Williams", "SG,SF,PF,C", "Oklahoma City Thunder"),
    @("Shai Gilgeous-Alexander", "PG", "Oklahoma City Thunder"),
    @("Dennis Schröder", "PG", "Brooklyn Nets"),
    @("Kyrie Irving", "PG,SG", "Dallas Mavericks"),
    @("Jordan Poole", "PG,SG", "Washington Wizards"),
    @("Brandon Boston Jr.", "SG,SF", "New Orleans Pelicans"),
    @("RJ Barrett", "SF,PF", "Toronto Raptors"),
    @("Jimmy Butler", "SF,PF", "Miami Heat"),
    @("CJ McCollum", "PG,SG", "New Orleans Pelicans"),
    @("Joel Embiid", "C", "Philadelphia 76ers")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
